$p = $ppt.ActivePresentation

# --- Slide 4 ("Justificación"): update the justification paragraph text and
#     grow the textbox so the longer paragraph still fits (spAutoFit). ---
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$shp4.TextFrame.TextRange.Text = "La implementación de un sistema de gestión de minutas en el SENA permitirá optimizar la organización y el uso eficiente del tiempo en los diferentes ambientes institucionales, beneficiando de manera directa tanto a los aprendices como a los instructores. Asimismo, el desarrollo de una plataforma digital para el registro de asistencia contribuirá a mejorar el control y seguimiento de los procesos formativos. Adicionalmente, el sistema ofrecerá un espacio para registrar apuntes o novedades relacionadas con incidentes, lo que fortalecerá la comunicación y el manejo oportuno de situaciones relevantes dentro de la institución."
$shp4.Height = 201.1453

# --- Slide 5 ("Objetivo General"): replace the general-objective text (also
#     dropping the bold/black formatting override, plain run now) and nudge
#     the textbox height to match the re-flowed text. ---
$s5 = $p.Slides.Item(5)
$shp5 = $s5.Shapes.Item(2)
$shp5.TextFrame.TextRange.Text = "Diseñar e implementar un sistema digital para la gestión de minutas, el control de asistencia y el registro de incidentes en los ambientes del SENA, con el fin de optimizar la organización institucional, mejorar la comunicación entre actores y facilitar el seguimiento eficiente de las actividades formativas."
$shp5.TextFrame.TextRange.Font.Bold = $false
$shp5.TextFrame.TextRange.LanguageID = "es-MX"
$shp5.Height = 104.2078
